$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Archive")
$ws.Range("B9").Value = '入夜113-115年'
$ws.Range("C9").Value = '入夜113-115年'
$ws.Range("D9").Value = 'N.F.113-115年'
$ws.Range("E9").Value = '입야 113-115년'
$ws.Range("F9").Value = 'N.F. 113-115'
$ws.Range("B24").Value = '秘金'
$ws.Range("C24").Value = '秘金'
$ws.Range("D24").Value = 'アダマント'
$ws.Range("E24").Value = '아카니움'
$ws.Range("F24").Value = 'Aurum Secretum'
$ws.Range("H24").Value = '秘金是一种费沙专属的尖端技术产物，拥有强大的共鸣效应，结合辅助材料及装置可以制作不同的秘金产品。'
$ws.Range("I24").Value = '秘金是一種費沙專屬的尖端技術產物，擁有強大的共鳴效應，結合輔助材料及裝置可以製作不同的秘金產品。'
$ws.Range("J24").Value = 'アダマントはフィシャ独自の先端技術の産物で、強力な共鳴効果を持つ。補助材料や装置と組み合わせることで、様々なアダマント製品を製造できる。'
$ws.Range("K24").Value = '아카니움은 피셔 전용 첨단 기술의 산물로 강력한 공명 효과를 가졌으며, 보조 재료와 장치를 결합해 다양한 아카니움 제품을 만들 수 있다.'
$ws.Range("L24").Value = 'Aurum Secretum is an advanced technological product exclusive to Fraser, possessing powerful resonance effects. With the assistance of auxiliary materials and devices, it can be used to create various Aurum Secretum products.'
$ws.Range("B38").Value = '市议会'
$ws.Range("C38").Value = '市議會'
$ws.Range("D38").Value = '市議会'
$ws.Range("E38").Value = '시의회'
$ws.Range("F38").Value = 'City Council'
$ws.Range("H38").Value = '狄斯城最高权力机构，由狄斯城公民直选产生，设有多项监督与政策裁定职权，受上庭监督。'
$ws.Range("I38").Value = '狄斯城最高權力機構，由狄斯城公民直選產生，設有多項監督與政策裁定職權，受上庭監督。'
$ws.Range("J38").Value = 'ディスシティの最高権力機関であり、ディスシティ市民によって直接選出される。複数の監督権および政策決定権を有し、上庭の監督を受ける。'
$ws.Range("K38").Value = '디스시티 최고 권력 기관, 디스시티 시민이 직접 선출하며, 여러 감독 및 정책 책정 직권을 가졌으며, 낙원의 감독을 받는다.'
$ws.Range("L38").Value = 'DisCity''s highest governing body, directly elected by its citizens, that holds multiple supervisory and policy-making powers while remaining under the oversight of Paradeisos.'
$ws.Range("B39").Value = '新城三区'
$ws.Range("C39").Value = '新城三區'
$ws.Range("D39").Value = 'ニューシティエリア3'
$ws.Range("E39").Value = '신성 3구역'
$ws.Range("F39").Value = 'Eastside District 3'
$ws.Range("H39").Value = '新城三区为狄斯城的政治中心，是所有立法、行政与司法权力的集中地。'
$ws.Range("I39").Value = '新城三區為狄斯城的政治中心，是所有立法、行政與司法權力的集中地。'
$ws.Range("J39").Value = 'ニューシティエリア3はディスシティの政治の中心であり、あらゆる立法・行政・司法権が集中する場所である。'
$ws.Range("K39").Value = '신성 3구역은 디스시티의 정치 중심지로, 모든 입법, 행정, 사법 권력이 집중된 지역이다.'
$ws.Range("L39").Value = 'Eastside District 3 serves as DisCity''s political hub, where legislative, administrative, and judicial authority is centralized.'
$ws.Range("B40").Value = '万国会'
$ws.Range("C40").Value = '萬國會'
$ws.Range("D40").Value = '博覧会'
$ws.Range("E40").Value = '만국 박람회'
$ws.Range("F40").Value = 'The Great Expo'
$ws.Range("H40").Value = '万国会为入夜34年由狄斯发起的跨城邦交流峰会，每四年一届，旨在恢复人类互联网络、推动资源与科技共享。'
$ws.Range("I40").Value = '萬國會為入夜34年由狄斯發起的跨城邦交流峰會，每四年一屆，旨在恢復人類交流網路、推動資源與科技共享。'
$ws.Range("J40").Value = '万国博覧会は、N.F.34年にディスが発起した都市間交流サミットで、4年に一度開催される。人類のコミュニケーションネットワークの再構築、資源とテクノロジーの共有・促進を目的としている。'
$ws.Range("K40").Value = '만국 박람회는 입야 34년 디스가 주도해 개최한 도시국가 간의 교류를 위한 정상회담이다. 4년에 한 번 개최되며 인류 네트워크 회복, 자원 및 과학 기술의 공유 촉진을 목적으로 한다.'
$ws.Range("L40").Value = 'The Great Expo is an inter-city-state exchange summit initiated by DisCity in N.F. 34, held every four years with the goal of restoring human communication networks and fostering resource and technology sharing.'
$ws.Range("B51").Value = '费沙'
$ws.Range("C51").Value = '費沙'
$ws.Range("D51").Value = 'フィシャ'
$ws.Range("E51").Value = '피셔'
$ws.Range("F51").Value = 'Fraser'
$ws.Range("H51").Value = '费沙是狄斯之外的大陆贸易中枢，有许多小城邦依赖费沙或受到其直接控制。'
$ws.Range("I51").Value = '費沙是狄斯之外的大陸貿易中樞，有許多小城邦依賴費沙或受到其直接控制。'
$ws.Range("J51").Value = 'フィシャはディス外の大陸における貿易の中心であり、多くの小都市がフィシャに依存するか、直接的に支配されている。'
$ws.Range("K51").Value = '피셔는 디스를 제외한 대륙 무역의 중추로, 수많은 소형 도시국가가 피셔에 의존하거나 직접적인 통제를 받는다.'
$ws.Range("L51").Value = 'Fraser is the Outland trade hub outside of DisCity, with many smaller cities either dependent on Fraser or under its direct control.'

$ws = $wb.Worksheets.Item("Entry")
$ws.Range("B33").Value = '冬月冲击系列案件'
$ws.Range("C33").Value = '冬月衝擊系列案件'
$ws.Range("D33").Value = '冬期連続襲撃事件'
$ws.Range("E33").Value = '겨울철 충격 시리즈 사건'
$ws.Range("F33").Value = 'Election Season Assaults'
$ws.Range("G33").Value = '114年11月28日，以内阁刺杀案爆发为起点，围绕议长大选，一连串袭击案件在新城三区接连爆发。'
$ws.Range("H33").Value = '114年11月28日，以內閣刺殺案爆發為起點，圍繞議長大選，一連串襲擊案件在新城三區接連爆發。'
$ws.Range("I33").Value = '114年11月28日、内閣暗殺事件を皮切りに、議長選挙を巡る襲撃事件がニューシティエリア3で相次いで発生した。'
$ws.Range("J33").Value = '114년 11월 28일, 내각 암살 사건을 기점으로, 의장 선거를 둘러싼 일련의 습격 사건이 신성 3구역에서 연이어 발생했다.'
$ws.Range("K33").Value = 'On November 28, N.F. 114, starting with the outbreak of the cabinet assassination case, a series of attack incidents erupted in succession across Eastside District 3—all tied to the Speaker Election.'
$ws.Range("B34").Value = '议长大选'
$ws.Range("C34").Value = '議長大選'
$ws.Range("D34").Value = '議長選挙'
$ws.Range("E34").Value = '의장 선거'
$ws.Range("F34").Value = 'Speaker Election'
$ws.Range("G34").Value = '114年12月26日，第22届城邦议长大选在经历了诸多风波后成功召开，迈克斯·金当选议长。'
$ws.Range("H34").Value = '114年12月26日，第22屆城邦議長大選在經歷了諸多風波後成功召開，邁克斯·金當選議長。'
$ws.Range("I34").Value = '114年12月26日、数々の混乱を経て、第22回ディスシティ議長選挙が無事に開催され、マイクス・キンが議長に選出された。'
$ws.Range("J34").Value = '114년 12월 26일, 수많은 혼란 끝에 제22회 도시국가 의상 선거가 무사히 개최되었고, 맥스·킹이 의장으로 당선되었다.'
$ws.Range("K34").Value = 'On December 26, N.F. 114, the 22nd Speaker Election was successfully held after weathering numerous upheavals, with Max King winning the position.'
$ws.Range("B35").Value = '万国会'
$ws.Range("C35").Value = '萬國會'
$ws.Range("D35").Value = '博覧会'
$ws.Range("E35").Value = '만국 박람회'
$ws.Range("F35").Value = 'The Great Expo'
$ws.Range("G35").Value = '115年6月16日，第13届万国会于狄斯开幕，异邦使节纷纷前往狄斯。'
$ws.Range("H35").Value = '115年6月16日，第13屆萬國會於狄斯開幕，異邦使節紛紛前往狄斯。'
$ws.Range("I35").Value = '115年6月16日、第13回万国博覧会がディスで開幕し、異邦の使節たちが続々とディスを訪れた。'
$ws.Range("J35").Value = '115년 6월 16일, 제13회 만국 박람회가 디스에서 개막했고, 이방의 사절들이 잇따라 디스로 향했다.'
$ws.Range("K35").Value = 'On June 16, N.F. 115, the 13th Great Expo commenced in DisCity, drawing envoys from the Outlands to the city.'
$ws.Range("B112").Value = '秘金'
$ws.Range("C112").Value = '秘金'
$ws.Range("D112").Value = 'アダマント'
$ws.Range("E112").Value = '아카니움'
$ws.Range("F112").Value = 'Aurum Secretum'
$ws.Range("G112").Value = '秘金是一种费沙专属的尖端技术产物，拥有强大的共鸣效应，结合辅助材料及装置可以制作不同的秘金产品。'
$ws.Range("H112").Value = '秘金是一種費沙專屬的尖端技術產物，擁有強大的共鳴效應，結合輔助材料及裝置可以製作不同的秘金產品。'
$ws.Range("I112").Value = 'アダマントはフィシャ独自の先端技術の産物で、強力な共鳴効果を持つ。補助材料や装置と組み合わせることで、様々なアダマント製品を製造できる。'
$ws.Range("J112").Value = '아카니움은 피셔 전용 첨단 기술의 산물로 강력한 공명 효과를 가졌으며, 보조 재료와 장치를 결합해 다양한 아카니움 제품을 만들 수 있다.'
$ws.Range("K112").Value = 'Aurum Secretum is an advanced technological product exclusive to Fraser, possessing powerful resonance effects. With the assistance of auxiliary materials and devices, it can be used to create various Aurum Secretum products.'
$ws.Range("B113").Value = '特性'
$ws.Range("C113").Value = '特性'
$ws.Range("D113").Value = '特性'
$ws.Range("E113").Value = '특성'
$ws.Range("F113").Value = 'Characteristics'
$ws.Range("G113").Value = '秘金是一种具备共鸣效应与放大效应的金属物质，其具备多种激发形态，在常态下会呈现出自然的几何结构。
基于共鸣和放大效应，秘金与其他元素如异方晶、狂厄产物、灰烬结晶等进行耦合激发后，会引发不同的干涉效果。'
$ws.Range("H113").Value = '秘金是一種具備共鳴效應與放大效應的金屬物質，其具備多種激發形態，在常態下會呈現出自然的幾何結構。
基於共鳴和放大效應，秘金與其他元素如異方晶、狂厄產物、灰燼結晶等進行耦合激發後，會引發不同的干涉效果。'
$ws.Range("I113").Value = 'アダマントとは、共鳴効果と増幅効果を併せ持つ金属物質であり、複数の活性形態を有する。通常状態では、自然な幾何学的構造を呈する。
共鳴および増幅効果に基づき、アダマントはイレギュラーキューブや狂瞳の産物、灰の結晶などの他元素と結合して活性化されることで、様々な干渉効果を引き起こす。'
$ws.Range("J113").Value = '아카니움은 공명 효과와 증폭 효과를 지닌 금속 물질로, 여러 가지 활성화 형태를 가지고 있으며, 평상시에는 자연스러운 기하 구조를 띤다.
공명과 증폭 효과로 인해 아카니움은 큐브, 변이 산물, 잿더미 결정 등 다른 원소와 결합해 활성화될 경우 서로 다른 간섭 효과를 유발하게 된다.'
$ws.Range("K113").Value = 'Aurum Secretum is a metallic substance that exhibits both resonance and amplification properties. It has multiple forms and naturally exhibits geometric structures in its inert state.
Due to its resonant and amplifying effects, when Aurum Secretum undergoes coupled activation with other elements, such as Hypercube, Mania byproducts, and Ash Crystal, it generates diverse interference phenomena.'
$ws.Range("B114").Value = '来源'
$ws.Range("C114").Value = '來源'
$ws.Range("D114").Value = '由来'
$ws.Range("E114").Value = '출처'
$ws.Range("F114").Value = 'Source'
$ws.Range("G114").Value = '秘金的原材料直接源自灰烬使徒的碎片，借助地底提供的特殊符文与炼金术式，对原材料进行多阶段灭活、稀释、合金化与最终提纯便可以获得秘金产物。'
$ws.Range("H114").Value = '秘金的原材料直接源自灰燼使徒的碎片，透過地底提供的特殊符文與煉金術式，對原材料進行多階段滅活、稀釋、合金化與最終提純便可以獲得秘金產物。'
$ws.Range("I114").Value = 'アダマントの原材料は、灰の使徒の破片に由来する。地底から提供される特殊な符文と錬金術式を用い、原材料に対して多段階の不活性化、希釈、合金化を行い、最終的に精製することで、アダマント産物を得られる。'
$ws.Range("J114").Value = '아카니움 원재료의 직접적인 공급원은 잿더미 사도 조각이다. 지하가 제공하는 특수 룬과 연금 주술의 힘을 빌려 원재료에 여러 단계의 비활성화, 희석, 합금화와 최종 정련을 진행하면 아카니움 산물을 얻을 수 있다.'
$ws.Range("K114").Value = 'The raw materials used for Aurum Secretum are sourced directly from Apostle of Ashes fragments. Through specialized runes and alchemical sorcery supplied by the Underground, the raw materials undergo multi-stage processing—including deactivation, dilution, alloying, and final purification—to yield the finished Aurum Secretum products.'
$ws.Range("B115").Value = '作用'
$ws.Range("C115").Value = '作用'
$ws.Range("D115").Value = '効果'
$ws.Range("E115").Value = '기능'
$ws.Range("F115").Value = 'Effect'
$ws.Range("G115").Value = '1.高纯度秘金与灰烬结晶结合形成的共鸣操控装置，使特定人员或禁闭者在短时间内可以与灰烬使徒的投影同调，干预其行动逻辑，甚至直接操纵灰烬使徒。
2.秘金制造的防护装置可以一定程度上防护与抵消狂厄，但秘金防护装置同样需要消耗异方晶。因其本质只是共鸣与放大异方晶的防护效果，无法在根本上避免对异方晶的依赖。'
$ws.Range("H115").Value = '1.高純度秘金與灰燼結晶結合形成的共鳴操控裝置，使特定人員或禁閉者在短時間內可以與灰燼使徒的投影同調，干預其行動邏輯，甚至直接操縱灰燼使徒。
2.秘金製造的防護裝置可以一定程度上防護與抵消狂厄，但秘金防護裝置同樣需要消耗異方晶。因其本質只是共鳴與放大異方晶的防護效果，無法在根本上避免對異方晶的依賴。'
$ws.Range("I115").Value = '1.高純度アダマントと灰の結晶を結合させた共鳴制御装置は、特定の人員またはコンビクトが短時間だけ灰の使徒の投影と同期し、その行動ロジックに干渉、さらには灰の使徒自体を直接操作することを可能にする。
2.アダマント製の防護装置は、狂瞳病に対する防護や相殺効果をある程度有しているが、同時にイレギュラーキューブの消費を必要とする。この本質は、イレギュラーキューブの防護効果を共鳴・増幅させるものであり、イレギュラーキューブへの依存を根本的に解消することはできない。'
$ws.Range("J115").Value = '1. 고순도 아카니움과 잿더미 결정이 결합해 만들어진 공명 조종 장치로, 특정 인원 또는 수감자가 짧은 시간 동안 잿더미 사도의 투영과 동조할 수 있게 만들어, 잿더미 사도의 행동 논리를 간섭하거나 잿더미 사도를 직접적으로 조종할 수 있다.
2. 아카니움으로 제조한 방어장치는 변이를 어느 정도 막아내거나 상쇄할 수 있지만, 아카니움 방어장치 역시 큐브를 소모해야 한다. 아카니움의 본질이 그저 큐브의 방호 효과를 공명하고 증폭하는 것에 불과하기 때문에, 근본적으로 큐브에 대한 의존을 벗어날 수는 없다.'
$ws.Range("K115").Value = '1. Resonance control devices, crafted by integrating high-purity Aurum Secretum with Ash Crystal, permit designated personnel or Sinners to briefly synchronize with Apostle of Ashes projections. This synchronization enables them to disrupt the Apostles'' behavioral logic or even exert direct control over them.
2. Protective devices constructed from Aurum Secretum offer partial defense and mitigation against Mania. However, these Aurum Secretum-based protections also demand Hypercube consumption. As they primarily function by resonating with and amplifying Hypercube''s protective properties, they remain incapable of fully eliminating reliance on Hypercube.'
$ws.Range("B179").Value = '市议会'
$ws.Range("C179").Value = '市議會'
$ws.Range("D179").Value = '市議会'
$ws.Range("E179").Value = '시의회'
$ws.Range("F179").Value = 'City Council'
$ws.Range("G179").Value = '狄斯城最高权力机构，由狄斯城公民直选产生，设有多项监督与政策裁定职权，受上庭监督。'
$ws.Range("H179").Value = '狄斯城最高權力機構，由狄斯城公民直選產生，設有多項監督與政策裁定職權，受上庭監督。'
$ws.Range("I179").Value = 'ディスシティの最高権力機関であり、ディスシティ市民によって直接選出される。複数の監督権および政策決定権を有し、上庭の監督を受ける。'
$ws.Range("J179").Value = '디스시티 최고 권력 기관, 디스시티 시민이 직접 선출하며, 여러 감독 및 정책 책정 직권을 가졌으며, 낙원의 감독을 받는다.'
$ws.Range("K179").Value = 'DisCity''s highest governing body, directly elected by its citizens, that holds multiple supervisory and policy-making powers while remaining under the oversight of Paradeisos.'
$ws.Range("B180").Value = '职能'
$ws.Range("C180").Value = '職能'
$ws.Range("D180").Value = '職務'
$ws.Range("E180").Value = '역할'
$ws.Range("F180").Value = 'Function'
$ws.Range("G180").Value = '市议会负责制定城市法规、预算审查与政策监督，拥有动议、审议、冻结和罢免内阁成员的权力。
其下设多个委员会，处理包括财政、治安、科技、灾害应对等事务，并可对涉及公民权利、公共资源分配与跨城邦事务的决策进行最终裁定。
议会对以下主要内阁部门拥有直接监督权与问责权：财政局、经贸局、治安总局、能源与工业局、外邦事务局与涅槃办事处。各部门须定期向议会提交运作报告，接受质询与预算审议，重要人事任免须经议会批准。'
$ws.Range("H180").Value = '市議會負責制定城市法規、預算審查與政策監督，擁有動議、審議、凍結和罷免內閣成員的權力。
其下設多個委員會，處理包括財政、治安、科技、災害應對等事務，並可對涉及公民權利、公共資源配置與跨城邦事務的決策進行最終裁定。
議會對以下主要內閣部門擁有直接監督權與問責權：財政局、經貿局、治安總局、能源與工業局、外邦事務局與涅槃辦事處。各部門須定期向議會提交運作報告，接受質詢與預算審議，重要人事任免須經議會批准。'
$ws.Range("I180").Value = '市議会は都市法規の制定、予算審査、政策監督を担当し、内閣メンバーに対する動議、審議、凍結、さらには罷免の権限を有している。
その下には複数の委員会が設置されており、財政、治安、科学技術、災害対応などの業務を扱うほか、市民の権利、公的資源の分配、都市間の事案に関わる決定について最終的な裁定を下すことができる。
議会は以下の主要な内閣部門に対して直接的な監督権および責任を問う権限を有している：財政局、経済貿易局、治安総局、エネルギー産業局、外邦事務局、涅槃事務局。各部門は定期的に議会へ運営報告を提出し、質疑や予算審査を受けなければならない。また、重要な人事の任免には議会の承認が必要である。'
$ws.Range("J180").Value = '시의회는 도시 법규 제정, 예산 심사 및 정책 감독을 담당하며, 발의, 심의, 동결 및 내각 구성원 파면 등 권력을 가진다.
산하에 설립된 여러 위원회는 재정, 치안, 과학, 재해 대응 등 업무를 맡으며, 시민의 권력, 공공 자원 분배 및 도시국가 간 사무 정책에 관한 최종 결정을 내릴 수 있다.
의회는 다음 주요 내각 부서에 대해 직접적인 감독 및 문책 권한을 갖는다: 재정국, 경제 무역국, 치안 총국, 에너지 산업국, 외무국, 니르바나 사무국. 각 부서는 반드시 정기적으로 의회에 운영 보고서를 제출하고, 질의와 예산 심사를 받아야 하며, 중요 인사의 선임 및 파면은 반드시 의회의 허가를 거쳐야 한다.'
$ws.Range("K180").Value = 'The city council is responsible for establishing municipal regulations, reviewing budgets, and overseeing policies, with the authority to propose, deliberate, freeze, and remove cabinet members.
It operates multiple committees that manage affairs such as finance, public security, technology, disaster response, and others, and holds final decision-making authority on matters involving civil rights, public resource allocation, and inter-city-state relations.
The council exercises direct oversight and accountability over the following major cabinet departments: the Treasury Department, Bureau of Commerce and Trade, Public Security Bureau, Bureau of Energy and Industry, Outland Affairs Bureau, and Nirvana Office. Each department must regularly submit operational reports to the council, undergo questioning and budget reviews, and secure council approval for key personnel appointments and dismissals.'
$ws.Range("B181").Value = '议长'
$ws.Range("C181").Value = '議長'
$ws.Range("D181").Value = '議長'
$ws.Range("E181").Value = '의장'
$ws.Range("F181").Value = 'Speaker'
$ws.Range("G181").Value = '议长为狄斯城市议会实质上的最高民选官员，具备行政调度与议会协调双重职责。其通过议长委员会协调市议会各项职能，同时对接内阁政府，是政策发布、公共事务调解与危机管理的核心枢纽人物。
议长不得兼任议员，其任期为四年，可连任一次。'
$ws.Range("H181").Value = '議長為狄斯城市議會實質上的最高民選官員，具備行政調度與議會協調雙重職責。其透過議長委員會協調市議會各項職能，同時對接內閣政府，是政策發布、公共事務調解與危機管理的核心樞紐人物。
議長不得兼任議員，其任期為四年，可連任一次。'
$ws.Range("I181").Value = '議長はディス市議会における実質的な最高民選官であり、行政指揮と議会調整の両方の役割を担う。議長は議長委員会を通じて市議会の各機能を調整し、同時に内閣政府と連携して、政策発表、公共事務の調整、危機管理における中枢的役割を担う。
議長は議員を兼任することはできず、任期は4年、1回まで再任可能。'
$ws.Range("J181").Value = '의장은 디스시티 시의회의 실질적인 최고 민선 관료로서, 행정 조율 및 의회 조정이라는 두 가지 역할을 겸한다. 의장은 의장 위원회를 통해 시의회의 각종 기능을 조율하고, 동시에 내각 정부와도 연결되어 정책 발표, 공공 업무 조정 및 위기관리의 핵심 허브 역할을 수행한다.
의장은 의원을 겸임할 수 없으며, 임기는 4년이고, 1회 연임할 수 있다.'
$ws.Range("K181").Value = 'The Speaker of the City Council serves as DisCity''s highest elected official in practice, carrying dual responsibilities for both administrative and legislative coordination. Through the Speaker''s Committee, the Speaker oversees the city council''s functions while also acting as the primary liaison with the cabinet government. In this role, the Speaker serves as the central hub for policy announcements, public affairs mediation, and crisis management.
The Speaker cannot simultaneously hold office as a councilor, serves a four-year term, and may be reelected for one additional consecutive term.'
$ws.Range("B182").Value = '议长选举'
$ws.Range("C182").Value = '議長選舉'
$ws.Range("D182").Value = '議長選挙'
$ws.Range("E182").Value = '의장 선거'
$ws.Range("F182").Value = 'Speaker Election'
$ws.Range("G182").Value = '议长不由议会内部选出，而由全城选民直接票选产生。每位登记选民拥有一票议长票，与议员选举同步进行。
1. 所有年满18周岁、在狄斯城连续居住超过五年、并拥有狄斯公民身份的居民，均享有议长选举投票权。
2. 现役FAC成员在服役期间将暂时冻结议长选举投票权。投票权将在服役期满、正式退出服役状态后自动解冻。'
$ws.Range("H182").Value = '議長不由議會內部選出，而由全城選民直接票選產生。每位登記選民擁有一票議長票，與議員選舉同步進行。
1. 所有年滿18周歲、在狄斯城連續居住超過五年、並擁有狄斯公民身分的居民，均享有議長選舉投票權。
2. 現役FAC成員在服役期間將暫時凍結議長選舉投票權。投票權將在服役期滿、正式退出服役狀態後自動解凍。'
$ws.Range("I182").Value = '議長は議会内部で選出されるのではなく、都市の有権者による直接投票で選出される。登録された全ての有権者は議長票を1票有し、議長選挙は議員選挙と同時に実施される。
1.ディスシティに継続して5年以上居住し、ディス市民権を有する18歳以上の全ての住民は、議長選挙の投票権を持つ。
2.現役のFACメンバーは任期中、議長選挙の投票権を一時的に凍結される。任期終了後、正式に退役した時点で投票権が自動的に凍結解除される。'
$ws.Range("J182").Value = '의장은 의회 내부에서 선출되지 않고, 도시 내 모든 유권자의 직접 투표로 선출된다. 등록된 모든 유권자는 의장 투표권 1표를 가지며, 의원 선거와 동시에 진행된다.
1. 만 18세 이상이고 디스시티에서 연속 5년 넘게 거주하고 디스 시민 신분을 보유한 모든 주민은 의장 선거 투표권을 갖는다.
2. 복무 중인 현역 FAC 대원은 의장 선거 투표권이 정지된다. 투표권은 복무 만료 후 공식적으로 퇴역할 시 자동으로 회복된다.'
$ws.Range("K182").Value = 'The Speaker of the City Council is not elected internally by the council but is instead chosen through a direct citywide popular vote. Each registered voter casts one ballot for Speaker, with the election held concurrently with councilor elections.
1. Eligibility is granted to any resident aged 18 or above, who has maintained continuous residence in DisCity for over five years and possesses DisCitizenship.
2. Active FAC members will have their voting rights temporarily suspended for the duration of their service. These rights will be reinstated upon completion of service and formal discharge from active duty.'
$ws.Range("B183").Value = '114年特别法案'
$ws.Range("C183").Value = '114年特別法案'
$ws.Range("D183").Value = '114年特別法案'
$ws.Range("E183").Value = '114년 특별 법안'
$ws.Range("F183").Value = 'N.F. 114 Special Act'
$ws.Range("G183").Value = '为进一步明确选民资格中“居住年限”条件中对狄斯城范围的认定，本法案对特定区域的居住年限计算规则作出如下补充：
对于新纳入狄斯城邦管辖范围的区域，其居民的“居住年限”自该区域正式纳入狄斯市议会管理之日起开始计算。
若某区域曾阶段性脱离议会管理，其居民的居住年限计算将以该区域重新纳入管理之日起为准。'
$ws.Range("H183").Value = '為進一步明確選民資格中「居住年限」條件中對狄斯城範圍的認定，本法案對特定區域的居住年限計算規則作出如下補充：
對於新納入狄斯城邦管轄範圍的區域，其居民的「居住年限」自該區域正式納入狄斯市議會管理之日起開始計算。
若某區域曾階段性脫離議會管理，其居民的居住年限計算將以該區域重新納入管理之日起為準。'
$ws.Range("I183").Value = '有権者資格における「居住年数」の条件に関連し、ディスシティの範囲認定をより明確にするため、本法案は特定エリアにおける居住年数の算定ルールに以下の補足を加える：
ディスの管轄範囲に新たに編入されたエリアについては、その住民の「居住年数」は、当該エリアが正式にディス市議会の管理下に置かれた日を起点として算定される。
あるエリアが一時的に議会の管理から外れていた場合、その住民の居住年数の算定は、当該エリアが再び管理下に編入された日を起点とする。'
$ws.Range("J183").Value = '유권자 자격 중 ''거주 기간'' 조건에서 디스시티 범위의 인정을 더 명확히 하기 위하여 본 법안에서 특정 구역의 거주 기간 계산 규칙을 다음과 같이 추가한다.
디스시티 관할 범위에 새롭게 추가된 구역에 관하여 해당 주민의 ''거주 기간''은 해당 구역이 정식으로 디스 시의회 관리하에 편입된 날부터 계산한다.
특정 구역이 단계적으로 의회 관리에서 벗어났을 경우, 해당 주민의 거주 기간 계산은 해당 구역이 다시 관리에 포함된 날을 기준으로 계산한다.'
$ws.Range("K183").Value = 'To further clarify the definition of DisCity''s jurisdiction regarding the residency period criterion for voter eligibility, this Act establishes the following supplementary rules for calculating residency periods in specific regions:
For regions recently included within DisCity''s jurisdiction, residents'' duration of residency shall be calculated starting from the official date on which the region comes under City Council administration.
For regions that were once temporarily separated from council administration, residents'' duration of residency shall be calculated starting from the date of reintegration into administration.'
$ws.Range("B184").Value = '新城三区'
$ws.Range("C184").Value = '新城三區'
$ws.Range("D184").Value = 'ニューシティエリア3'
$ws.Range("E184").Value = '신성 3구역'
$ws.Range("F184").Value = 'Eastside District 3'
$ws.Range("G184").Value = '新城三区为狄斯城的政治中心，是所有立法、行政与司法权力的集中地。该区域围绕中央公园规划，设有议长官邸、狄斯纪念堂、新城治安总局、狄斯档案馆等。'
$ws.Range("H184").Value = '新城三區為狄斯城的政治中心，是所有立法、行政與司法權力的集中地。該區域圍繞中央公園規劃，設有議長官邸、狄斯紀念堂、新城治安總局、狄斯檔案館等。'
$ws.Range("I184").Value = 'ニューシティエリア3はディスシティの政治の中心であり、あらゆる立法・行政・司法権が集中する場所である。このエリアは中央公園を中心に整備されており、議長官邸、ディス記念堂、ニューシティ治安総局、ディス公文書館などの施設が設置されている。'
$ws.Range("J184").Value = '신성 3구역은 디스시티의 정치 중심지로, 모든 입법, 행정, 사법 권력이 집중된 지역이다. 해당 구역은 중앙공원을 중심으로 설계되었으며, 의장 관저, 디스 기념관, 신성 치안 총국, 디스 기록관 등이 위치해 있다.'
$ws.Range("K184").Value = 'Eastside District 3 serves as DisCity''s political hub, where legislative, administrative, and judicial authority is centralized. Designed around Central Park, the district accommodates key facilities such as the Speaker''s Residence, DisCity Memorial Hall, Eastside Public Security Bureau, and the DisCity Archives.'
$ws.Range("B185").Value = '议长官邸'
$ws.Range("C185").Value = '議長官邸'
$ws.Range("D185").Value = '議長官邸'
$ws.Range("E185").Value = '의장 관저'
$ws.Range("F185").Value = 'Speaker''s Residence'
$ws.Range("G185").Value = '议长官邸位于新城三区中央公园北侧中部，是市议会运行建筑群的核心，兼具议会官邸与议会大厅的功能，东西两翼分层设有各委员会及内阁会议室，各行政机关办事处以及城邦应急指挥中心。'
$ws.Range("H185").Value = '議長官邸位於新城三區中央公園北側中部，是市議會運行建築群的核心，兼具議會官邸與議會大廳的功能，東西兩翼分層設有各委員會及內閣會議室，各行政機關辦事處以及城邦應急指揮中心。'
$ws.Range("I185").Value = '議長官邸はニューシティエリア3の中央公園北側中央部に位置し、市議会運営棟の中核を成す建物だ。議会官邸と議会ホールの機能を兼ね備えており、東西には階層ごとに各委員会室および内閣会議室、各行政機関の事務室、さらに都市の緊急指揮センターが設けられている。'
$ws.Range("J185").Value = '의장 관저는 신성 3구역 중앙공원 북측 중앙에 있으며, 시의회 운영 건물들의 핵심으로, 의회 관저 및 의회 로비 기능을 겸한다. 동서 양쪽 날개에는 각 위원회 및 내각 회의실과, 각 행정 기관 사무실과 도시국가 비상 지휘 센터가 마련되어 있다.'
$ws.Range("K185").Value = 'The Speaker''s Residence is situated in the central area on the north side of Central Park in Eastside District 3. It forms the heart of the city council''s operational complex, serving as both the Speaker''s official residence and the council chamber. The east and west wings span multiple floors, accommodating committee and cabinet meeting rooms, offices for administrative agencies, and the city-state emergency command center.'
$ws.Range("B186").Value = '冬月冲击'
$ws.Range("C186").Value = '冬月衝擊'
$ws.Range("D186").Value = '冬期襲撃'
$ws.Range("E186").Value = '겨울철 충격'
$ws.Range("F186").Value = 'December Surprises'
$ws.Range("G186").Value = '入夜114年冬月，议长大选前，新城三区内连续发生数起袭击事件。以内阁袭击案为起始，血色哀宴案件为高潮，连续不断的政治暴力严重冲击了该区域的治安环境。
在冬月冲击事件后，新城三区强化了安保等级，除强化狂厄监测机制之外还额外设置了严密的各类监控体系以应对异邦渗透。'
$ws.Range("H186").Value = '入夜114年冬月，議長大選前，新城三區內連續發生數起襲擊事件。以內閣襲擊案為起始，血色哀宴案件為高潮，連續不斷的政治暴力嚴重衝擊了該區域的治安環境。
在冬月衝擊事件後，新城三區強化了保全等級，除強化狂厄監測機制之外還額外設置了嚴密的各類監控體系以應對異邦滲透。'
$ws.Range("I186").Value = 'N.F.114年の冬、議長選挙を前にして、ニューシティエリア3で連続襲撃事件が発生した。内閣襲撃事件に始まり、血の哀宴事件に至るまでの相次ぐ政治的暴力は、当該エリアの治安環境に深刻な打撃を与えた。
冬期連続襲撃事件の後、ニューシティエリア3は警備体制を強化した。狂瞳病の監視を強固にしただけでなく、異邦からの侵入に備えて厳重な各種監視システムも追加で導入した。'
$ws.Range("J186").Value = '입야 114년 겨울, 의장 선거를 앞두고 신성 3구역에서 연속적인 습격 사건이 발생했다. 내각 습격 사건을 시작으로 블러디 파티 사건에서 절정에 이르는, 끊임없는 정치적 폭력 충돌이 해당 구역의 치안 환경에 심각한 타격을 주었다.
겨울철 충격 사건 이후, 신성 3구역은 보안 등급을 강화했으며, 변이 감시 메커니즘을 강화하는 것 외에도 추가로 여러 치밀한 감시 체계를 만들어 이방 침입에 대응했다.'
$ws.Range("K186").Value = 'During the winter months of N.F. 114, preceding the Speaker Election, a string of consecutive attacks struck Eastside District 3. Beginning with the cabinet assassination case and culminating in the Blood Memorial incident, the wave of relentless political violence severely destabilized the district''s security environment.
In the aftermath of the election season assaults, Eastside District 3 elevated its security protocols, deploying stricter Mania monitoring measures and instituting comprehensive surveillance systems to deter Outland infiltration.'
$ws.Range("B187").Value = '万国会'
$ws.Range("C187").Value = '萬國會'
$ws.Range("D187").Value = '博覧会'
$ws.Range("E187").Value = '만국 박람회'
$ws.Range("F187").Value = 'The Great Expo'
$ws.Range("G187").Value = '万国会为入夜34年由狄斯发起的跨城邦交流峰会，每四年一届，旨在恢复人类互联网络、推动资源与科技共享。历经战争与灾难波折，多届万国会停办。
入夜115年6月，在费沙的推动下，新城特意在西区奇兰广场附近的首届万国会旧址建立了水晶宫，召开了第13届万国会。'
$ws.Range("H187").Value = '萬國會為入夜34年由狄斯發起的跨城邦交流峰會，每四年一屆，旨在恢復人類交流網路、推動資源與科技共享。歷經戰爭與災難波折，多屆萬國會停辦。
入夜115年6月，在費沙的推動下，新城特意在西區奇蘭廣場附近的首屆萬國會舊址建立了水晶宮，召開了第13屆萬國會。'
$ws.Range("I187").Value = '万国博覧会は、N.F.34年にディスが発起した都市間交流サミットで、4年に一度開催される。人類のコミュニケーションネットワークの再構築、資源とテクノロジーの共有・促進を目的としている。戦争や災害による波乱で、博覧会は複数回、開催中止となった。
N.F.115年6月、フィシャの後押しのもと、ニューシティは特別に西区のチラン広場付近にある第1回万国博覧会の跡地にクリスタルパレスを建設し、第13回万国博覧会を開催した。'
$ws.Range("J187").Value = '만국 박람회는 입야 34년 디스가 주도해 개최한 도시국가 간의 교류를 위한 정상회담이다. 4년에 한 번 개최되며 인류 네트워크 회복, 자원 및 과학 기술의 공유 촉진을 목적으로 한다. 전쟁과 재난으로 인해 여러 차례 만국 박람회의 개최가 중단되기도 했다.
입야 115년 6월, 피셔의 주도 아래 신성은 서구 쉘란 플라자 근처, 첫 만국 박람회가 열렸던 곳에 특별히 수정궁을 건설하고 제13회 만국 박람회를 개최했다.'
$ws.Range("K187").Value = 'The Great Expo is an inter-city-state exchange summit initiated by DisCity in N.F. 34, held every four years with the goal of restoring human communication networks and fostering resource and technology sharing. However, due to the upheavals of wars and disasters, multiple sessions of the Great Expo were suspended.
In June, N.F. 115, driven by Fraser''s initiative, Eastside constructed Cristallia near the original site of the first Great Expo in West District''s Keylan Square area to host the 13th Great Expo.'
$ws.Range("B188").Value = '历史'
$ws.Range("C188").Value = '歷史'
$ws.Range("D188").Value = '歴史'
$ws.Range("E188").Value = '역사'
$ws.Range("F188").Value = 'History'
$ws.Range("G188").Value = '奇兰广场建于入夜33年，以纪念奇兰博士，次年，第一届万国会于此召开，共有40余城邦受邀，新生的狄斯城向全世界展示了异方晶体系，开启了入夜时代各城邦的合作可能。'
$ws.Range("H188").Value = '奇蘭廣場建於入夜33年，以紀念奇蘭博士，次年，第一屆萬國會於此召開，共有40餘城邦受邀，新生的狄斯城向全世界展示了異方晶體系，開啟了入夜時代各城邦的合作可能。'
$ws.Range("I188").Value = 'チラン広場はN.F.33年、チラン博士を称えて建設された。翌年、ここで第1回万国博覧会が開催され、40を超える都市が招待された。新たに生まれたディスシティは、イレギュラーキューブシステムを世界に披露し、ナイトフォール時代における各都市間の協力の可能性を切り開いた。'
$ws.Range("J188").Value = '쉘란 플라자는 입야 33년, 쉘란 박사를 기념하기 위해 건설되었다. 이듬해 제1회 만국 박람회가 이곳에서 개최되었으며, 총 40여 개의 도시국가가 초청되었다. 새롭게 탄생한 디스시티는 전 세계에 큐브 체계를 선보이며, 입야 시대 각 도시국가 간의 협력 가능성을 만들었다.'
$ws.Range("K188").Value = 'Keylan Square was constructed in N.F. 33 in honor of Dr. Keylan. The following year, it hosted the inaugural Great Expo, welcoming over 40 city-states. The fledgling DisCity unveiled the Hypercube technology to the world, paving the way for collaboration among city-states in the Nightfall Era.'
$ws.Range("G189").Value = '自入夜60年代起，异邦战争导致万国会多次中断，仅有极少届作为临时的和平信号召开。
入夜106年，第12届万国会召开，虽然经历变故与挫折，但最终还是取得了积极的意义，次年各城邦成功签订全面停战。'
$ws.Range("H189").Value = '自入夜60年代起，異邦戰爭導致萬國會多次中斷，僅有極少屆作為臨時的和平訊號召開。
入夜106年，第12屆萬國會召開，雖然經歷變故與挫折，但最終還是取得了積極的意義，次年各城邦成功簽訂全面停戰。'
$ws.Range("I189").Value = 'N.F.60年代以降、異邦戦争の影響により、博覧会はたびたび中断され、一時的な平和の象徴としてわずか数回のみ開催された。
N.F.106年には第12回万国博覧会が開催され、数々の混乱や挫折を経験しながらも、最終的には一定の成果を上げ、翌年には各都市間で全面停戦協定が締結された。'
$ws.Range("J189").Value = '입야 60년대부터 시작된 이방 전쟁으로 인해 만국 박람회는 수차례 중단되었으며, 단 몇 차례만이 임시적인 평화의 신호로서 간헐적으로 개최되었다.
입야 106년, 제12회 만국 박람회가 열렸고, 혼란과 좌절을 겪었음에도 끝내 의미 있는 성과를 이루어냈다. 이듬해 각 도시국가는 전면 정전 협정에 합의하게 되었다.'
$ws.Range("K189").Value = 'Beginning in the N.F. 60s, the Outland wars repeatedly disrupted the Great Expo, with only a handful of sessions held as temporary gestures of peace.
In N.F. 106, the 12th Great Expo was convened. Despite facing setbacks and challenges, it ultimately yielded positive outcomes, and the following year, all city-states successfully signed a comprehensive ceasefire agreement.'
$ws.Range("B190").Value = '水晶宫'
$ws.Range("C190").Value = '水晶宮'
$ws.Range("D190").Value = 'クリスタルパレス'
$ws.Range("E190").Value = '수정궁'
$ws.Range("F190").Value = 'Cristallia'
$ws.Range("G190").Value = '水晶宫分为主馆与分馆，虽然两大场馆分立两侧，各自提供不同的公共展品，但中间有无数连廊沟通，象征各个城邦之间的互通有无，完美照应本届万国会“再次链接世界”的主题。
除公共展区外，水晶宫高层设置有内部展区，专向各城邦使团及新城政商代表开放。'
$ws.Range("H190").Value = '水晶宮分為主館與分館，雖然兩大場館分立兩側，各自提供不同的公共展品，但中間有無數連廊溝通，象徵各個城邦之間的互通有無，完美照應本屆萬國會「再次連結世界」的主題。
除公共展區外，水晶宮高層設置有內部展區，專向各城邦使團及新城政商代表開放。'
$ws.Range("I190").Value = 'クリスタルパレスは本館と別館に分かれている。二つの大展示館は左右に分かれてそれぞれ異なる公共展示が行われているものの、その間には無数の連絡回廊が設けられており、各都市間の相互交流と連携を象徴している。これは、今回の博覧会のテーマである「再び世界を繋ぐ」に完璧に合致している。
公共展示エリアのほか、クリスタルパレスの上層階には内部展示エリアが設けられており、各都市の使節団やニューシティの政財界代表に限定して公開されている。'
$ws.Range("J190").Value = '수정궁은 본관과 별관으로 나뉘며, 두 전시관은 양측에 위치해 각각 다른 공공 전시품을 선보이지만, 중간에 각 도시국가 간의 소통을 상징하는 무수한 통로가 연결되어 있어, 이번 만국 박람회의 주제인 ''다시 세계를 연결한다''를 완벽하게 표현하고 있다.
공공 전시 구역 외에도 수정궁 고층에는 각 도시국가 사절단과 신성의 정치 및 상업 대표에게만 개방되는 내부 전시 구역이 마련되어 있다.'
$ws.Range("K190").Value = 'Cristallia is divided into a main hall and an auxiliary hall. While the two facilities stand on opposite sides, each hosting distinct public exhibitions, countless connecting corridors unite them—symbolizing the mutual exchange and communication between city-states and perfectly reflecting this Great Expo''s theme: "Be One Again."
In addition to the public exhibition areas, Cristallia''s upper floors house private exhibition spaces, accessible exclusively to delegations from various city-states and Eastside''s political and business representatives.'
$ws.Range("B191").Value = '闭幕式袭击'
$ws.Range("C191").Value = '閉幕式襲擊'
$ws.Range("D191").Value = '閉幕式の襲撃'
$ws.Range("E191").Value = '폐막식 습격'
$ws.Range("F191").Value = 'Closing Ceremony Attack'
$ws.Range("G191").Value = '入夜115年6月25日，万国会闭幕式上，异邦极端分子对展会发动了亡命袭击，但在预先设置的安保措施下，袭击者成功被FAC控制。'
$ws.Range("H191").Value = '入夜115年6月25日，萬國會閉幕式上，異邦極端分子對展會發動了亡命襲擊，但在預先設置的保全措施下，襲擊者成功被FAC控制。'
$ws.Range("I191").Value = 'N.F.115年6月25日、万国博覧会の閉幕式において、異邦の過激派による自爆襲撃が行われた。しかし、事前に講じられていた警備措置により、襲撃者はFACによって無事制圧された。'
$ws.Range("J191").Value = '입야 115년 6월 25일, 만국 박람회 폐막식 도중, 이방의 극단 세력이 전시회에 테러 습격을 감행했지만, 미리 준비된 보안 조치 덕에 FAC가 이를 신속히 제압하였다.'
$ws.Range("K191").Value = 'On June 25, N.F. 115, during the Great Expo''s closing ceremony, Outland extremists launched a desperate attack on the exhibition, but the assailants were swiftly subdued by the FAC due to preemptive security measures.'
$ws.Range("B229").Value = '费沙'
$ws.Range("C229").Value = '費沙'
$ws.Range("D229").Value = 'フィシャ'
$ws.Range("E229").Value = '피셔'
$ws.Range("F229").Value = 'Fraser'
$ws.Range("G229").Value = '费沙位于大陆西北边缘，是世界上除了狄斯以外的最大城邦，依靠强大的贸易、矿业及工业能力在城邦联盟之中居于核心地位。'
$ws.Range("H229").Value = '費沙位於大陸西北邊緣，是世界上除了狄斯以外的最大城邦，依靠強大的貿易、礦業及工業能力在城邦聯盟之中居於核心地位。'
$ws.Range("I229").Value = 'フィシャは大陸の北西端に位置し、ディスを除けば世界最大の都市である。強力な貿易、鉱業、工業能力で都市連盟の中核を担っている。'
$ws.Range("J229").Value = '피셔는 대륙 북서쪽 변두리에 있는, 디스 다음으로 세계에서 가장 큰 도시국가이다. 강력한 무역, 광산업, 공업 능력을 기반으로 도시국가 연맹에서 핵심적인 지위를 차지하고 있다.'
$ws.Range("K229").Value = 'Fraser is situated at the northwestern edge of the continent and is the world''s largest city-state after DisCity. It maintains a dominant role within the Farland Alliance due to its formidable trade, mining, and industrial strength.'
$ws.Range("B230").Value = '城市信息'
$ws.Range("C230").Value = '城市資訊'
$ws.Range("D230").Value = '都市情報'
$ws.Range("E230").Value = '도시 정보'
$ws.Range("F230").Value = 'City Info'
$ws.Range("G230").Value = '费沙是狄斯之外的大陆贸易中枢，有许多小城邦依赖费沙或受到其直接控制。
1.产业
依靠旧时代的经济底蕴，费沙保持了较为完整的科技及工业基础，其针对砂海特制的贸易砂船终日不停，穿梭在各个城邦之间。
除此之外，费沙的金融业与奢侈品行业也非常发达。
2.人口
费沙核心城区人口只有约百万级，主城区外分布着大量附属卫星城及小型村落，大量的边缘劳动人口与灾后难民都聚居于此。据上一次统计数据显示，整个费沙体系的人口数量已经在N.F.110年突破1200万，由于缺乏官方统计，实际数据仍存在一定争议。
3.经济
费沙控制着大陆西北多个重要矿脉，并通过贸易协定与砂海专行的机制获得大量异邦输入财富，也是难得的可以与狄斯城进行其他各类产品贸易的城邦。'
$ws.Range("H230").Value = '費沙是狄斯之外的大陸貿易中樞，有許多小城邦依賴費沙或受到其直接控制。
1.產業
依靠舊時代的經濟底蘊，費沙保持了較為完整的科技及工業基礎，其針對砂海特製的貿易砂船終日不停，穿梭在各個城邦之間。
除此之外，費沙的金融業與奢侈品行業也非常發達。
2.人口
費沙核心城區人口只有約百萬級，主城區外分佈著大量附屬衛星城及小型村落，大量的邊緣勞動人口與災後難民都聚居于此。據上一次統計資料顯示，整個費沙體系的人口數量已經在N.F.110年突破1200萬，由於缺乏官方統計，實際數據仍存在一定爭議。
3.經濟
費沙控制著大陸西北多個重要礦脈，並透過貿易協定與砂海專行的機制獲得大量異邦輸入財富，也是難得的可以與狄斯城進行其他各類產品貿易的城邦。'
$ws.Range("I230").Value = 'フィシャはディス外の大陸における貿易の中心であり、多くの小都市がフィシャに依存するか、直接的に支配されている。
1.産業
フィシャは旧時代の経済基盤に支えられ、比較的充実した科学技術と工業基礎を維持している。砂の海専用に設計された貿易用の砂航船は、昼夜を問わず稼働し、各都市間を行き来している。
その他にも、フィシャは金融業や嗜好品産業においても非常に発展している。
2.人口
フィシャの中心エリアの人口はおよそ100万人程度だが、中心エリア外には多数の付属衛星都市や小さな村があり、そこに周縁労働者や災害を受けた難民が居住している。直近の統計によれば、フィシャ体系の全体の人口はN.F.110年に1200万人を突破したが、公式な統計がないため、実際の数値については議論の余地がある。
3.経済
フィシャは大陸北西部にある複数の重要な鉱脈を掌握しており、貿易協定や砂の海独自の制度を通じて、多くの異邦から金を流入させている。また、ディスシティとその他各種製品の貿易を行うことができる、貴重な都市の一つでもある。'
$ws.Range("J230").Value = '피셔는 디스를 제외한 대륙 무역의 중추로, 수많은 소형 도시 국가들이 피셔에 의존하거나 직접적인 통제를 받는다.
1. 산업
구시대의 경제 기반 덕에 피셔는 비교적 온전한 기술 및 공업 기반을 유지했으며, 화이트 샌드에 특화된 피셔의 무역 사막 함선은 쉬지 않고 각 도시국가 사이를 오가고 있다.
이 외에 피셔의 금융업과 사치품 산업 또한 매우 발달했다.
2. 인구
피셔 핵심 구역의 인구는 약 백만 명 규모이지만, 주요 도심 구역 밖에는 수많은 부속 위성도시와 소형 마을이 산재해 있으며, 대량의 비정규 노동자 인구와 재난을 겪은 난민들이 이곳에 거주하고 있다. 가장 최근 통계에 따르면, 피셔 체계의 인구수는 이미 N.F.110년에 1,200만 명을 돌파했다. 그러나 공식적인 통계가 없으므로 실제 수치와는 다를 수 있다.
3. 경제
피셔는 대륙 북서쪽의 여러 주요 광맥을 통제하고 있으며, 무역 협정과 화이트 샌드 전용 메커니즘을 통해 이방으로부터 수많은 부를 흡수했다. 또한 피셔는 디스시티와 다양한 상품 무역이 가능한 몇 안 되는 도시국가 중 하나이기도 하다.'
$ws.Range("K230").Value = 'Fraser is the Outland trade hub outside of DisCity, with many smaller cities either dependent on Fraser or under its direct control.
1. Industry
Leveraging its old-world economic foundations, Fraser has preserved relatively intact technological and industrial infrastructure. Its specialized sand ships, designed for traversing WhiteSands, operate continuously, shuttling between various city-states.
Additionally, Fraser boasts a highly developed financial sector and luxury goods industry.
2. Population
Fraser''s core urban area has a population of only about one million, while numerous affiliated satellite towns and small settlements sprawl beyond the main city. These peripheral areas house large numbers of marginalized workers and post-calamity refugees. According to the latest census data, the entire Fraser system''s population exceeded 12 million in N.F. 110, though the actual figures remain disputed due to the lack of official statistics.
3. Economy
Fraser controls multiple key mineral veins in the northwest and amasses substantial wealth from the Outlands through trade agreements and WhiteSands'' exclusive transit mechanisms. It is also one of the few city-states capable of conducting diversified trade with DisCity.'

